# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" header suffixes to the respective format
# version names ("_FV2404" / "_FV2410"), turns the used range into a
# native Excel Table (ListObject), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 "_old" -> "_FV2404", L1:U1 "_new" -> "_FV2410") ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# --- 2. Freeze the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into a table ---
$rng = $ws.Range("A1:U73")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
